$d = $word.ActiveDocument

$d.Content.Find.Execute("2022-12-20 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2022-12-21 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("74-4=", $true, $false, $false, $false, $false, $true, 1, $false, "65-44=", 2) | Out-Null
$d.Content.Find.Execute("12+38=", $true, $false, $false, $false, $false, $true, 1, $false, "90-49=", 2) | Out-Null
$d.Content.Find.Execute("92-63=", $true, $false, $false, $false, $false, $true, 1, $false, "90-23=", 2) | Out-Null
$d.Content.Find.Execute("85-40=", $true, $false, $false, $false, $false, $true, 1, $false, "79-51=", 2) | Out-Null
$d.Content.Find.Execute("78+7=", $true, $false, $false, $false, $false, $true, 1, $false, "43+2=", 2) | Out-Null
$d.Content.Find.Execute("88-14=", $true, $false, $false, $false, $false, $true, 1, $false, "48-25=", 2) | Out-Null
$d.Content.Find.Execute("51-8=", $true, $false, $false, $false, $false, $true, 1, $false, "72-10=", 2) | Out-Null
$d.Content.Find.Execute("6+29=", $true, $false, $false, $false, $false, $true, 1, $false, "65-24=", 2) | Out-Null
$d.Content.Find.Execute("22-2=", $true, $false, $false, $false, $false, $true, 1, $false, "58+22=", 2) | Out-Null
$d.Content.Find.Execute("10+26=", $true, $false, $false, $false, $false, $true, 1, $false, "97-97=", 2) | Out-Null
$d.Content.Find.Execute("7+10=", $true, $false, $false, $false, $false, $true, 1, $false, "69+18=", 2) | Out-Null
$d.Content.Find.Execute("55+40=", $true, $false, $false, $false, $false, $true, 1, $false, "75-41=", 2) | Out-Null
$d.Content.Find.Execute("61-11=", $true, $false, $false, $false, $false, $true, 1, $false, "49-48=", 2) | Out-Null
$d.Content.Find.Execute("68+11=", $true, $false, $false, $false, $false, $true, 1, $false, "66+16=", 2) | Out-Null
$d.Content.Find.Execute("10+45=", $true, $false, $false, $false, $false, $true, 1, $false, "46-9=", 2) | Out-Null
$d.Content.Find.Execute("93-90=", $true, $false, $false, $false, $false, $true, 1, $false, "5+93=", 2) | Out-Null
$d.Content.Find.Execute("37+4=", $true, $false, $false, $false, $false, $true, 1, $false, "27-23=", 2) | Out-Null
$d.Content.Find.Execute("44-9=", $true, $false, $false, $false, $false, $true, 1, $false, "54+45=", 2) | Out-Null
$d.Content.Find.Execute("74-54=", $true, $false, $false, $false, $false, $true, 1, $false, "58-15=", 2) | Out-Null
$d.Content.Find.Execute("68-48=", $true, $false, $false, $false, $false, $true, 1, $false, "84-18=", 2) | Out-Null
$d.Content.Find.Execute("0+50=", $true, $false, $false, $false, $false, $true, 1, $false, "58-6=", 2) | Out-Null
$d.Content.Find.Execute("41+35=", $true, $false, $false, $false, $false, $true, 1, $false, "65-53=", 2) | Out-Null
$d.Content.Find.Execute("65-36=", $true, $false, $false, $false, $false, $true, 1, $false, "34+27=", 2) | Out-Null
$d.Content.Find.Execute("48-9=", $true, $false, $false, $false, $false, $true, 1, $false, "7+12=", 2) | Out-Null
$d.Content.Find.Execute("49+27=", $true, $false, $false, $false, $false, $true, 1, $false, "65-42=", 2) | Out-Null
$d.Content.Find.Execute("64-7=", $true, $false, $false, $false, $false, $true, 1, $false, "75-32=", 2) | Out-Null
$d.Content.Find.Execute("5+57=", $true, $false, $false, $false, $false, $true, 1, $false, "83-8=", 2) | Out-Null
$d.Content.Find.Execute("37+61=", $true, $false, $false, $false, $false, $true, 1, $false, "90-80=", 2) | Out-Null
$d.Content.Find.Execute("53+9=", $true, $false, $false, $false, $false, $true, 1, $false, "97-59=", 2) | Out-Null
$d.Content.Find.Execute("91-10=", $true, $false, $false, $false, $false, $true, 1, $false, "32+21=", 2) | Out-Null
$d.Content.Find.Execute("56-42=", $true, $false, $false, $false, $false, $true, 1, $false, "59-30=", 2) | Out-Null
$d.Content.Find.Execute("76-26=", $true, $false, $false, $false, $false, $true, 1, $false, "23+31=", 2) | Out-Null
$d.Content.Find.Execute("85-38=", $true, $false, $false, $false, $false, $true, 1, $false, "90+7=", 2) | Out-Null
$d.Content.Find.Execute("39+35=", $true, $false, $false, $false, $false, $true, 1, $false, "29-28=", 2) | Out-Null
$d.Content.Find.Execute("37+42=", $true, $false, $false, $false, $false, $true, 1, $false, "89-12=", 2) | Out-Null
$d.Content.Find.Execute("49+49=", $true, $false, $false, $false, $false, $true, 1, $false, "21+45=", 2) | Out-Null
$d.Content.Find.Execute("6+42=", $true, $false, $false, $false, $false, $true, 1, $false, "86+5=", 2) | Out-Null
$d.Content.Find.Execute("38+55=", $true, $false, $false, $false, $false, $true, 1, $false, "14+9=", 2) | Out-Null
$d.Content.Find.Execute("78-64=", $true, $false, $false, $false, $false, $true, 1, $false, "61+17=", 2) | Out-Null
$d.Content.Find.Execute("24+9=", $true, $false, $false, $false, $false, $true, 1, $false, "16+41=", 2) | Out-Null
$d.Content.Find.Execute("15+60=", $true, $false, $false, $false, $false, $true, 1, $false, "3-2=", 2) | Out-Null
$d.Content.Find.Execute("26-20=", $true, $false, $false, $false, $false, $true, 1, $false, "6+57=", 2) | Out-Null
$d.Content.Find.Execute("31+48=", $true, $false, $false, $false, $false, $true, 1, $false, "57+24=", 2) | Out-Null
$d.Content.Find.Execute("86-45=", $true, $false, $false, $false, $false, $true, 1, $false, "63-2=", 2) | Out-Null
$d.Content.Find.Execute("81-7=", $true, $false, $false, $false, $false, $true, 1, $false, "30+2=", 2) | Out-Null
$d.Content.Find.Execute("33+62=", $true, $false, $false, $false, $false, $true, 1, $false, "10+52=", 2) | Out-Null
$d.Content.Find.Execute("11+62=", $true, $false, $false, $false, $false, $true, 1, $false, "8+20=", 2) | Out-Null
$d.Content.Find.Execute("70-69=", $true, $false, $false, $false, $false, $true, 1, $false, "97-77=", 2) | Out-Null
$d.Content.Find.Execute("22-3=", $true, $false, $false, $false, $false, $true, 1, $false, "62-37=", 2) | Out-Null
$d.Content.Find.Execute("99-2=", $true, $false, $false, $false, $false, $true, 1, $false, "33+48=", 2) | Out-Null
$d.Content.Find.Execute("61+9=", $true, $false, $false, $false, $false, $true, 1, $false, "61-57=", 2) | Out-Null
$d.Content.Find.Execute("96-32=", $true, $false, $false, $false, $false, $true, 1, $false, "86-66=", 2) | Out-Null
$d.Content.Find.Execute("14+57=", $true, $false, $false, $false, $false, $true, 1, $false, "38+26=", 2) | Out-Null
$d.Content.Find.Execute("83-45=", $true, $false, $false, $false, $false, $true, 1, $false, "96+0=", 2) | Out-Null
$d.Content.Find.Execute("75-35=", $true, $false, $false, $false, $false, $true, 1, $false, "79-46=", 2) | Out-Null
$d.Content.Find.Execute("82-14=", $true, $false, $false, $false, $false, $true, 1, $false, "73-58=", 2) | Out-Null
$d.Content.Find.Execute("66+2=", $true, $false, $false, $false, $false, $true, 1, $false, "21+71=", 2) | Out-Null
$d.Content.Find.Execute("59-15=", $true, $false, $false, $false, $false, $true, 1, $false, "60-2=", 2) | Out-Null
$d.Content.Find.Execute("38+11=", $true, $false, $false, $false, $false, $true, 1, $false, "35-14=", 2) | Out-Null
$d.Content.Find.Execute("2+89=", $true, $false, $false, $false, $false, $true, 1, $false, "39+1=", 2) | Out-Null
$d.Content.Find.Execute("37-5=", $true, $false, $false, $false, $false, $true, 1, $false, "57-33=", 2) | Out-Null
$d.Content.Find.Execute("42+25=", $true, $false, $false, $false, $false, $true, 1, $false, "48+51=", 2) | Out-Null
$d.Content.Find.Execute("87-73=", $true, $false, $false, $false, $false, $true, 1, $false, "75-53=", 2) | Out-Null
$d.Content.Find.Execute("62-38=", $true, $false, $false, $false, $false, $true, 1, $false, "53+23=", 2) | Out-Null
$d.Content.Find.Execute("82-30=", $true, $false, $false, $false, $false, $true, 1, $false, "92+6=", 2) | Out-Null
$d.Content.Find.Execute("22+38=", $true, $false, $false, $false, $false, $true, 1, $false, "9+66=", 2) | Out-Null
$d.Content.Find.Execute("65-14=", $true, $false, $false, $false, $false, $true, 1, $false, "95-63=", 2) | Out-Null
$d.Content.Find.Execute("85-70=", $true, $false, $false, $false, $false, $true, 1, $false, "77-31=", 2) | Out-Null
$d.Content.Find.Execute("61+20=", $true, $false, $false, $false, $false, $true, 1, $false, "63+33=", 2) | Out-Null
$d.Content.Find.Execute("35+38=", $true, $false, $false, $false, $false, $true, 1, $false, "43+41=", 2) | Out-Null
$d.Content.Find.Execute("66+12=", $true, $false, $false, $false, $false, $true, 1, $false, "91-51=", 2) | Out-Null
$d.Content.Find.Execute("24-3=", $true, $false, $false, $false, $false, $true, 1, $false, "77-42=", 2) | Out-Null
$d.Content.Find.Execute("8+15=", $true, $false, $false, $false, $false, $true, 1, $false, "63-33=", 2) | Out-Null
$d.Content.Find.Execute("16+11=", $true, $false, $false, $false, $false, $true, 1, $false, "16-11=", 2) | Out-Null
$d.Content.Find.Execute("51-13=", $true, $false, $false, $false, $false, $true, 1, $false, "92-1=", 2) | Out-Null
$d.Content.Find.Execute("41-19=", $true, $false, $false, $false, $false, $true, 1, $false, "51+29=", 2) | Out-Null
$d.Content.Find.Execute("58-27=", $true, $false, $false, $false, $false, $true, 1, $false, "99-6=", 2) | Out-Null
$d.Content.Find.Execute("49+22=", $true, $false, $false, $false, $false, $true, 1, $false, "75+22=", 2) | Out-Null
$d.Content.Find.Execute("62-52=", $true, $false, $false, $false, $false, $true, 1, $false, "80-71=", 2) | Out-Null
$d.Content.Find.Execute("66-33=", $true, $false, $false, $false, $false, $true, 1, $false, "34+57=", 2) | Out-Null
$d.Content.Find.Execute("31-27=", $true, $false, $false, $false, $false, $true, 1, $false, "95-70=", 2) | Out-Null
$d.Content.Find.Execute("57+34=", $true, $false, $false, $false, $false, $true, 1, $false, "51+37=", 2) | Out-Null
$d.Content.Find.Execute("43-36=", $true, $false, $false, $false, $false, $true, 1, $false, "0+70=", 2) | Out-Null
$d.Content.Find.Execute("70-45=", $true, $false, $false, $false, $false, $true, 1, $false, "18+77=", 2) | Out-Null
$d.Content.Find.Execute("56-49=", $true, $false, $false, $false, $false, $true, 1, $false, "93-74=", 2) | Out-Null
$d.Content.Find.Execute("11+72=", $true, $false, $false, $false, $false, $true, 1, $false, "26+17=", 2) | Out-Null
$d.Content.Find.Execute("57+25=", $true, $false, $false, $false, $false, $true, 1, $false, "45-30=", 2) | Out-Null
$d.Content.Find.Execute("4+88=", $true, $false, $false, $false, $false, $true, 1, $false, "46+29=", 2) | Out-Null
$d.Content.Find.Execute("22+27=", $true, $false, $false, $false, $false, $true, 1, $false, "30+58=", 2) | Out-Null
$d.Content.Find.Execute("88-64=", $true, $false, $false, $false, $false, $true, 1, $false, "73-19=", 2) | Out-Null
$d.Content.Find.Execute("47+25=", $true, $false, $false, $false, $false, $true, 1, $false, "33+16=", 2) | Out-Null
$d.Content.Find.Execute("25+50=", $true, $false, $false, $false, $false, $true, 1, $false, "70-60=", 2) | Out-Null
$d.Content.Find.Execute("28-15=", $true, $false, $false, $false, $false, $true, 1, $false, "29+31=", 2) | Out-Null
$d.Content.Find.Execute("84+15=", $true, $false, $false, $false, $false, $true, 1, $false, "15-8=", 2) | Out-Null
$d.Content.Find.Execute("10+24=", $true, $false, $false, $false, $false, $true, 1, $false, "60-29=", 2) | Out-Null
$d.Content.Find.Execute("97-43=", $true, $false, $false, $false, $false, $true, 1, $false, "77+2=", 2) | Out-Null
$d.Content.Find.Execute("37+17=", $true, $false, $false, $false, $false, $true, 1, $false, "2+17=", 2) | Out-Null
$d.Content.Find.Execute("86-72=", $true, $false, $false, $false, $false, $true, 1, $false, "25+58=", 2) | Out-Null
$d.Content.Find.Execute("12+47=", $true, $false, $false, $false, $false, $true, 1, $false, "57-42=", 2) | Out-Null
$d.Content.Find.Execute("30-0=", $true, $false, $false, $false, $false, $true, 1, $false, "7+64=", 2) | Out-Null
